# Split the bolded "**Discuss with Michael...**" run on the product-page
# bullet into three runs: the leading " **", the new task description, and
# the trailing "**" -- replacing the old sentence with the new one while
# keeping the bold markers.

$d = $word.ActiveDocument

$oldPhrase = " **Discuss with Michael how to implement this as it might already exist**"
$newMiddle = "Add a quantity indicator to the product page and build a function for the  enter button, next to the quantity value that hides the bag it button if a value higher than the quantity in stock Is input"

function Escape-Xml([string]$text) {
    return $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
}

# Locate the run/sentence that needs to be split.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target phrase to replace."
}
$hitStart = $rng.Start
$hitEnd = $rng.End

# Find the paragraph that contains the hit so we can rebuild it whole
# (this keeps the paragraph's own formatting / numbering / style intact).
$paragraph = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $hitStart -and $p.Range.End -ge $hitEnd) {
        $paragraph = $p
        break
    }
}
if ($null -eq $paragraph) {
    throw "Could not locate the paragraph containing the target phrase."
}

$paraRange = $paragraph.Range
$fullText = $paraRange.Text
# Paragraph.Range.Text includes the trailing paragraph mark; drop it before
# comparing / rebuilding.
$paraMark = [string][char]13
if ($fullText.EndsWith($paraMark)) {
    $fullText = $fullText.Substring(0, $fullText.Length - 1)
}
if (-not $fullText.EndsWith($oldPhrase)) {
    throw "Unexpected paragraph content; refusing to edit."
}
$prefixText = $fullText.Substring(0, $fullText.Length - $oldPhrase.Length)

# Pull the existing <w:pPr> (style/numbering/run formatting default) out of
# the paragraph's own OOXML so the rebuilt paragraph keeps it unchanged.
$paraOoxml = $paraRange.WordOpenXML
$pPr = ""
if ($paraOoxml -match '(?s)<w:p\b[^>]*>\s*(<w:pPr>.*?</w:pPr>)') {
    $pPr = $matches[1]
}

$prefixRun = ""
if ($prefixText.Length -gt 0) {
    $prefixRun = '<w:r><w:t xml:space="preserve">' + (Escape-Xml $prefixText) + '</w:t></w:r>'
}

$newRuns = ''
$newRuns += '<w:r><w:t xml:space="preserve"> **</w:t></w:r>'
$newRuns += '<w:r><w:t>' + (Escape-Xml $newMiddle) + '</w:t></w:r>'
$newRuns += '<w:r><w:t>**</w:t></w:r>'

$bodyInner = $pPr + $prefixRun + $newRuns

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $bodyInner + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$null = $paraRange.InsertXML($packageXml)
